# Refresh the cryptocurrency price/volume snapshot (rows 2-51).
# Price (col D) and Volume(1h) (col E) are stored as plain text in the sheet,
# so numeric-looking prices are written with a leading apostrophe to keep them
# as text, then the style is reset to Normal so no visual "quote prefix" marker
# or number formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.551.15"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.912.23"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'247.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.19%  "
$ws.Range("D6").Value = "'0.658"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.63%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "'42.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Value = "'0.347"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.16%  "
$ws.Range("D10").Value = "'49.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.78%  "
$ws.Range("D11").Value = "'0.0718"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.36%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "2.187.93"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "'12.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.29%  "
$ws.Range("D15").Value = "'0.703"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.83%  "
$ws.Range("D16").Value = "1.904.84"
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").Value = "'4.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.85%  "
$ws.Range("D18").Value = "35.544.13"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").Value = "'72.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "0.0₃0834"
$ws.Range("E20").Value = "  +4.79%  "
$ws.Range("D21").Value = "'245.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("D22").Value = "'12.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.72%  "
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +17.13%  "
$ws.Range("D27").Value = "'171.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "'8.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.17%  "
$ws.Range("D29").Value = "'18.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("D30").Value = "'0.128"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.29%  "
$ws.Range("E31").Value = "  +4.39%  "
$ws.Range("D32").Value = "'0.967"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +23.01%  "
$ws.Range("D33").Value = "'0.0572"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("D34").Value = "'4.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.25%  "
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").Value = "'1.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.37%  "
$ws.Range("D37").Value = "'2.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("D40").Value = "'0.0206"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").Value = "'92.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "'0.0639"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +18.00%  "
$ws.Range("D43").Value = "'15.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.17%  "
$ws.Range("D44").Value = "1.349.22"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D46").Value = "'47.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +38.30%  "
$ws.Range("D47").Value = "'12.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "'2.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").Value = "'2.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").Value = "'6.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").Value = "2.097.59"
$ws.Range("E51").Value = "  +2.83%  "
